# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp text update
$ws.Range("A1").Value = "Datos actualizados a 22 de Mayo de 2020 a las 00:35"

# Estados Unidos (row 4) - updated daily totals
$ws.Range("B4").Value = 1617856
$ws.Range("C4").Value = 25133
$ws.Range("D4").Value = 381676
$ws.Range("E4").Value = 1139985
$ws.Range("G4").Value = 1259
$ws.Range("H4").Value = 96195

# Canada (row 17) - updated daily totals
$ws.Range("B17").Value = 81312
$ws.Range("C17").Value = 1170
$ws.Range("D17").Value = 41696
$ws.Range("E17").Value = 33467
$ws.Range("G17").Value = 118
$ws.Range("H17").Value = 6149

# Re-sort of tied countries (Montserrat / Groenlandia / Seychelles) moves their
# rows; data travels with the new row order.
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 0

$ws.Range("A210").Value = "Seychelles"

$ws.Range("A211").Value = "Montserrat"
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1

# Re-sort of tied countries (Sahara Occidental / Bonaire, San Eustaquio y Saba)
$ws.Range("A214").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("A215").Value = "Sahara Occidental"
